$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "INSERT INTO ..." statements for column O (rows 2-11).
# The Zip4 placeholder changed from '' to NULL (bug fix for UT pass-through).
$inserts = @(
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('1','1700 W. Washington','','','Phoenix','AZ','85007', NULL,'1')",
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('2','820 N. French St.','','','Wilmington','DE','19801', NULL,'0')",
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('3','700 W. Jefferson','','','Boise','ID','83720',NULL,'1')",
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('4','700 Capitol Ave.','','','Frankfurt','KY','40601',NULL,'1')",
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('5','75 Constitution Ave.','','','St. Paul','MN','55155',NULL,'0')",
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('6','101 N. Carson Street','','','Carson City','NV','89701',NULL,'1')",
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('7','116 West Jones St.','','','Raleigh','NC','27603',NULL,'0')",
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('8','600 E. Boulevard Avenue','','','Bismark','ND','58505',NULL,'0')",
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('9','225 Main Capitol','','','Harrisburg','PA','17120',NULL,'0')",
    "INSERT INTO [VolTeer].[Vol].[tblVolAddress] (AddrID,AddrLine1,AddrLine2,AddrLine3,City,St,Zip,Zip4,ActiveFlg) VALUES('10','109 State St.','','','Montpelier','VT','56409',NULL,'0')"
)

for ($i = 0; $i -lt $inserts.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 15)
    $cell.Value = $inserts[$i]
    $cell.Style = "Normal"
}

# Reflect final selection/viewport seen in the saved workbook.
$ws.Range("O11").Select()
try {
    $excel.ActiveWindow.ScrollColumn = 13
} catch {
    # Scroll position isn't critical; ignore if unsupported.
}
